$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.794.08'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '1.763.27'
$ws.Range("E3").Value = '  -2.70%  '
$ws.Range("D4").Value = "'1.006"
$ws.Range("D5").Value = "'339.13"
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").Value = "'0.3757"
$ws.Range("E7").Value = '  -4.06%  '
$ws.Range("D8").Value = "'0.3359"
$ws.Range("D9").Value = "'45.74"
$ws.Range("E9").Value = '  -5.26%  '
$ws.Range("D10").Value = "'1.132"
$ws.Range("E10").Value = '  -5.66%  '
$ws.Range("D11").Value = "'0.07269"
$ws.Range("E11").Value = '  -3.81%  '
$ws.Range("D12").Value = "'22.75"
$ws.Range("E12").Value = '  +2.77%  '
$ws.Range("D13").Value = "'1.001"
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").Value = "'6.234"
$ws.Range("E14").Value = '  -4.38%  '
$ws.Range("D15").Value = "'7.218"
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("D16").Value = '1.766.04'
$ws.Range("E16").Value = '  -2.57%  '
$ws.Range("D17").Value = "'0.00001053"
$ws.Range("E17").Value = '  -4.78%  '
$ws.Range("D18").Value = "'0.06595"
$ws.Range("E18").Value = '  -1.47%  '
$ws.Range("D19").Value = "'80.95"
$ws.Range("E19").Value = '  -4.84%  '
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("D21").Value = "'17.11"
$ws.Range("D22").Value = "'6.329"
$ws.Range("E22").Value = '  -3.61%  '
$ws.Range("D23").Value = '27.877.75'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = "'11.78"
$ws.Range("E24").Value = '  -8.35%  '
$ws.Range("D25").Value = "'2.377"
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("D26").Value = "'1.489"
$ws.Range("E26").Value = '  +0.80%  '
$ws.Range("D27").Value = "'153.39"
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("D28").Value = "'20.01"
$ws.Range("E28").Value = '  -6.17%  '
$ws.Range("D29").Value = "'2.340"
$ws.Range("E29").Value = '  -8.36%  '
$ws.Range("D30").Value = '1.967.33'
$ws.Range("E30").Value = '  -2.57%  '
$ws.Range("D31").Value = "'131.64"
$ws.Range("E31").Value = '  -2.80%  '
$ws.Range("D32").Value = "'4.028"
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").Value = "'5.895"
$ws.Range("E33").Value = '  -3.53%  '
$ws.Range("D34").Value = "'0.08733"
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("D35").Value = "'12.37"
$ws.Range("E35").Value = '  -6.73%  '
$ws.Range("D36").Value = "'0.02359"
$ws.Range("E36").Value = '  -2.46%  '
$ws.Range("D37").Value = "'0.6716"
$ws.Range("E37").Value = '  -2.86%  '
$ws.Range("D38").Value = "'0.06243"
$ws.Range("E38").Value = '  -4.29%  '
$ws.Range("D39").Value = "'5.184"
$ws.Range("E39").Value = '  -5.98%  '
$ws.Range("D40").Value = "'0.2118"
$ws.Range("E40").Value = '  -4.74%  '
$ws.Range("D41").Value = "'1.219"
$ws.Range("E41").Value = '  -3.70%  '
$ws.Range("D42").Value = "'1.464"
$ws.Range("E42").Value = '  -9.38%  '
$ws.Range("D43").Value = "'8.064"
$ws.Range("E43").Value = '  -5.57%  '
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = '  +0.36%  '
$ws.Range("D45").Value = "'13.92"
$ws.Range("E45").Value = '  -5.43%  '
$ws.Range("D46").Value = "'0.6107"
$ws.Range("E46").Value = '  -6.32%  '
$ws.Range("D47").Value = "'3.843"
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("D48").Value = "'131.10"
$ws.Range("E48").Value = '  -1.04%  '
$ws.Range("D49").Value = "'2.027"
$ws.Range("E49").Value = '  -6.24%  '
$ws.Range("D50").Value = "'0.07269"
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("D51").Value = "'1.183"
$ws.Range("E51").Value = '  +1.54%  '
